$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '70.776.37'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.94%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.586.99'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.29%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '587.55'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.48%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '186.12'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.21%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.576.82'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.14%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.622'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.45%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.04%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +17.07%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.18%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '54.37'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.40%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +5.04%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.156.49'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.10%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '19.61'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.29%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '70.785.76'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.04%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.596.41'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.55%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '576.79'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +17.32%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.41'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.04%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.01%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.82%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '17.80'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -8.60%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.66'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +5.95%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.50%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '95.36'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -1.50%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.39'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.50%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.17%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.16'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -2.00%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '32.29'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +1.95%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.32'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -5.34%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '12.42'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +2.50%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '65.00'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.34%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.91%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.37'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +5.26%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '563.50'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -2.25%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +2.58%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '37.74'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -3.25%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.13%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0₃0795'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.03%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.391.36'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +6.63%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.135'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.76%  '
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.39'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -2.99%  '
$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.11'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.58%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.57'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -2.27%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -2.70%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.98'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -2.93%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.34%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +1.24%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.00'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.07%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.46'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -8.39%  '
